$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.757.53"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").Value = "2.271.24"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'303.77"
$ws.Range("E5").Value = "  +0.32%  "

# Row 6
$ws.Range("D6").Value = "'92.58"
$ws.Range("E6").Value = "  +0.96%  "

# Row 7
$ws.Range("E7").Value = "  +1.94%  "

# Row 8
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").Value = "'32.55"
$ws.Range("E10").Value = "  +1.20%  "

# Row 11
$ws.Range("D11").Value = "'53.37"
$ws.Range("E11").Value = "  -1.49%  "

# Row 12
$ws.Range("E12").Value = "  +0.28%  "

# Row 13
$ws.Range("E13").Value = "  -1.30%  "

# Row 14
$ws.Range("D14").Value = "'6.67"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15
$ws.Range("D15").Value = "2.624.57"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").Value = "'14.31"
$ws.Range("E16").Value = "  +1.32%  "

# Row 17
$ws.Range("D17").Value = "2.284.54"
$ws.Range("E17").Value = "  +1.29%  "

# Row 18
$ws.Range("D18").Value = "'0.777"
$ws.Range("E18").Value = "  +3.57%  "

# Row 19
$ws.Range("D19").Value = "41.682.82"
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("E20").Value = "  +1.69%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +0.14%  "

# Row 22
$ws.Range("E22").Value = "  +0.53%  "

# Row 23
$ws.Range("D23").Value = "'67.14"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").Value = "'243.31"
$ws.Range("E24").Value = "  +1.03%  "

# Row 25
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  +0.40%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.27%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'1.93"
$ws.Range("E27").Value = "  +3.53%  "

# Row 28
$ws.Range("D28").Value = "'23.96"
$ws.Range("E28").Value = "  +0.77%  "

# Row 29
$ws.Range("E29").Value = "  -1.36%  "

# Row 30
$ws.Range("E30").Value = "  -5.44%  "

# Row 31
$ws.Range("D31").Value = "'35.35"
$ws.Range("E31").Value = "  +5.34%  "

# Row 32
$ws.Range("D32").Value = "'160.51"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("E33").Value = "  +0.73%  "

# Row 34
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("D35").Value = "'0.0743"
$ws.Range("E35").Value = "  +0.92%  "

# Row 36
$ws.Range("E36").Value = "  -0.78%  "

# Row 37
$ws.Range("D37").Value = "'16.88"
$ws.Range("E37").Value = "  +1.18%  "

# Row 38
$ws.Range("E38").Value = "  +0.31%  "

# Row 39
$ws.Range("E39").Value = "  +0.88%  "

# Row 40
$ws.Range("E40").Value = "  +0.65%  "

# Row 41
$ws.Range("E41").Value = "  +0.40%  "

# Row 42
$ws.Range("D42").Value = "'3.91"
$ws.Range("E42").Value = "  -0.99%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.005.77"
$ws.Range("E43").Value = "  -2.77%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.36"
$ws.Range("E44").Value = "  -6.22%  "

# Row 45
$ws.Range("E45").Value = "  +1.99%  "

# Row 46
$ws.Range("D46").Value = "'10.32"
$ws.Range("E46").Value = "  +0.87%  "

# Row 47
$ws.Range("E47").Value = "  +3.31%  "

# Row 48
$ws.Range("E48").Value = "  -2.24%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'52.55"
$ws.Range("E49").Value = "  +3.08%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.51"
$ws.Range("E50").Value = "  -0.69%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.15"
$ws.Range("E51").Value = "  +0.92%  "
